$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3057.9038
$ws.Range("J17").Value = 3099.0393
$ws.Range("L17").Value = 9297.117899999999
$ws.Range("N17").Value = -9633.117899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 29573.027
$ws.Range("I64").Value = 145085.72
$ws.Range("J64").Value = 2620.0667
$ws.Range("K64").Value = 145085.72
$ws.Range("L64").Value = 2620.0667
$ws.Range("M64").Value = -144837.72
$ws.Range("N64").Value = -3116.0667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 29573.027
$ws.Range("I67").Value = 145085.72
$ws.Range("J67").Value = 2620.0667
$ws.Range("K67").Value = 145085.72
$ws.Range("L67").Value = 2620.0667
$ws.Range("M67").Value = -144227.72
$ws.Range("N67").Value = -4336.066699999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1336.6364
$ws.Range("I70").Value = 1166.6666
$ws.Range("J70").Value = 1540.6
$ws.Range("K70").Value = 3499.9998
$ws.Range("L70").Value = 4621.799999999999
$ws.Range("M70").Value = -3229.9998
$ws.Range("N70").Value = -5161.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1336.6364
$ws.Range("I73").Value = 1166.6666
$ws.Range("J73").Value = 1540.6
$ws.Range("K73").Value = 3499.9998
$ws.Range("L73").Value = 4621.799999999999
$ws.Range("M73").Value = -2563.9998
$ws.Range("N73").Value = -6493.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 20000
$ws.Range("J17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("N17").Value = -20346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 46396.5
$ws.Range("J113").Value = 46396.5
$ws.Range("L113").Value = 46396.5
$ws.Range("N113").Value = -55074.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 41994
$ws.Range("J114").Value = 41994
$ws.Range("L114").Value = 41994
$ws.Range("N114").Value = -50672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 54811
$ws.Range("J119").Value = 54811
$ws.Range("L119").Value = 54811
$ws.Range("N119").Value = -64487

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 31377.25
$ws.Range("J137").Value = 41600
$ws.Range("L137").Value = 41600
$ws.Range("N137").Value = -51800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 72999.60000000001
$ws.Range("J57").Value = 72999.60000000001
$ws.Range("L57").Value = 72999.60000000001
$ws.Range("N57").Value = -74439.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 43866.668
$ws.Range("J59").Value = 43866.668
$ws.Range("L59").Value = 43866.668
$ws.Range("N59").Value = -45560.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 48590
$ws.Range("J110").Value = 48590
$ws.Range("L110").Value = 48590
$ws.Range("N110").Value = -56770

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 55332.668
$ws.Range("J133").Value = 55332.668
$ws.Range("L133").Value = 55332.668
$ws.Range("N133").Value = -65452.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H136").Value = 72999.60000000001
$ws.Range("J136").Value = 72999.60000000001
$ws.Range("L136").Value = 72999.60000000001
$ws.Range("N136").Value = -83199.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 51000
$ws.Range("J139").Value = 51000
$ws.Range("L139").Value = 51000
$ws.Range("N139").Value = -61280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 20868.889
$ws.Range("I17").Value = 18955
$ws.Range("J17").Value = 22400
$ws.Range("K17").Value = 18955
$ws.Range("L17").Value = 22400
$ws.Range("M17").Value = -18781
$ws.Range("N17").Value = -22748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 42701.332
$ws.Range("J110").Value = 42701.332
$ws.Range("L110").Value = 42701.332
$ws.Range("N110").Value = -50881.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 44686
$ws.Range("J119").Value = 44686
$ws.Range("L119").Value = 44686
$ws.Range("N119").Value = -54362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 20438
$ws.Range("J137").Value = 20438
$ws.Range("L137").Value = 20438
$ws.Range("N137").Value = -30638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 32080
$ws.Range("J139").Value = 29600
$ws.Range("L139").Value = 29600
$ws.Range("N139").Value = -39880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 99999
$ws.Range("J141").Value = 99999
$ws.Range("L141").Value = 99999
$ws.Range("N141").Value = -110359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5633.8335
$ws.Range("I68").Value = 750
$ws.Range("J68").Value = 15401.5
$ws.Range("K68").Value = 2250
$ws.Range("L68").Value = 46204.5
$ws.Range("M68").Value = -1439
$ws.Range("N68").Value = -47826.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 5633.8335
$ws.Range("I71").Value = 750
$ws.Range("J71").Value = 15401.5
$ws.Range("K71").Value = 6750
$ws.Range("L71").Value = 138613.5
$ws.Range("M71").Value = -2694
$ws.Range("N71").Value = -146725.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 41499.5
$ws.Range("J114").Value = 41499.5
$ws.Range("L114").Value = 41499.5
$ws.Range("N114").Value = -50177.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1325.9565
$ws.Range("I122").Value = 1215.6316
$ws.Range("K122").Value = 3646.8948
$ws.Range("M122").Value = -1196.8948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 19938.46
$ws.Range("J137").Value = 19938.46
$ws.Range("L137").Value = 19938.46
$ws.Range("N137").Value = -30138.46

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 30268
$ws.Range("J139").Value = 30268
$ws.Range("L139").Value = 30268
$ws.Range("N139").Value = -40548

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 24450
$ws.Range("J64").Value = 24450
$ws.Range("L64").Value = 24450
$ws.Range("N64").Value = -24900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 24450
$ws.Range("J67").Value = 24450
$ws.Range("L67").Value = 24450
$ws.Range("N67").Value = -26010

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3599.8
$ws.Range("I68").Value = 2999.5
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2999.5
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -2250.5
$ws.Range("N68").Value = -5498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3599.8
$ws.Range("I71").Value = 2999.5
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 14997.5
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -11253.5
$ws.Range("N71").Value = -27488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 47408
$ws.Range("J119").Value = 47408
$ws.Range("L119").Value = 47408
$ws.Range("N119").Value = -57084

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 22888.562
$ws.Range("J133").Value = 22888.562
$ws.Range("L133").Value = 22888.562
$ws.Range("N133").Value = -27948.562

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 24305.455
$ws.Range("J137").Value = 24305.455
$ws.Range("L137").Value = 24305.455
$ws.Range("N137").Value = -34505.455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 43597.418
$ws.Range("J139").Value = 43597.418
$ws.Range("L139").Value = 43597.418
$ws.Range("N139").Value = -53877.418

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3522.2222
$ws.Range("J62").Value = 3522.2222
$ws.Range("L62").Value = 3522.2222
$ws.Range("N62").Value = -4770.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3522.2222
$ws.Range("J65").Value = 3522.2222
$ws.Range("L65").Value = 17611.111
$ws.Range("N65").Value = -23851.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 19908.846
$ws.Range("J139").Value = 19908.846
$ws.Range("L139").Value = 19908.846
$ws.Range("N139").Value = -30188.846
